$wb = $excel.ActiveWorkbook

$ws0 = $wb.Worksheets.Item("Layer0")
$ws0.Range("B2").Value = -0.1225721343096147
$ws0.Range("C2").Value = 0.5606186046271986
$ws0.Range("B3").Value = 0.7173803524541073
$ws0.Range("C3").Value = -1.404933774125337
$ws0.Range("B4").Value = -0.2909214003976496
$ws0.Range("C4").Value = -1.047722932434492

$ws1 = $wb.Worksheets.Item("Layer1")
$ws1.Range("B2").Value = -0.4842590857289767
$ws1.Range("C2").Value = 0.2506462288344503
$ws1.Range("B3").Value = 0.2734170130585224
$ws1.Range("C3").Value = -0.2967719627180544
$ws1.Range("B4").Value = -2.207334342030598
$ws1.Range("C4").Value = -0.2939936145181908
